# Apply the "add keys and reformat melatonin tables" edit to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Update the D13 text - append two additional lines to the existing note
$ws.Range("D13").Value = "Melatonin random-contolled n's & %'s -placebo AEs" + [char]10 + "Corrections and additions to melatonin data" + [char]10 + "Melatonin efficacy data"

# 2. Set the row height for row 13 to 57 to fit the new wrapped text
$ws.Rows.Item(13).RowHeight = 57

# 3. Update the hours value in E13 from 8 to 15
$ws.Range("E13").Value = 15

# Formulas in E15 (SUM) and E16 will recalculate automatically.
$excel.Calculate()
